# Regenerated data from tools
#
# The dataset_id=6 rows (rows 10-12: binary_understandability,
# correct_verif_questions, time_to_understand) were recomputed on each of
# the 5 "tool" sheets: num_snippets_warnings (F), num_warnings (G) and the
# derived correlation statistics (I kendalls_tau, J kendalls_p_value,
# K spearmans_rho, L spearmans_p_value) all changed. Two sheets also picked
# up incidental column-width tweaks (col I on all_tools, cols J/K on infer).

$wb = $excel.ActiveWorkbook

function Set-Row12Stats {
    param(
        $ws,
        [int]$row,
        $f,
        $g,
        [string]$i,
        [string]$j,
        [string]$k,
        [string]$l
    )

    if ($null -ne $f) { $ws.Cells.Item($row, 6).Value = $f }
    if ($null -ne $g) { $ws.Cells.Item($row, 7).Value = $g }
    $ws.Cells.Item($row, 9).Value  = [double]$i
    $ws.Cells.Item($row, 10).Value = [double]$j
    $ws.Cells.Item($row, 11).Value = [double]$k
    $ws.Cells.Item($row, 12).Value = [double]$l
}

# ---------------------------------------------------------------------
# all_tools
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("all_tools")

Set-Row12Stats $ws 10 48 752 "0.02919871531052824" "0.7792869073970613" "0.05326702354515665" "0.7133288232331647"
Set-Row12Stats $ws 11 48 752 "0.03091191513338553" "0.7563176058989638" "0.05189590320941469" "0.7204036820310951"
Set-Row12Stats $ws 12 48 752 "0.1841800495983816"  "0.06278475836348044" "0.2809257097381211" "0.04813238376383353"

# col I width 21.7109375 -> 20.7109375 (closest reachable value through the
# ColumnWidth property, which snaps to pixel increments)
$ws.Columns.Item(9).ColumnWidth = 19.75

# ---------------------------------------------------------------------
# checker_framework
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("checker_framework")

Set-Row12Stats $ws 10 23 68 "-0.1127953217534384" "0.3246303816462205" "-0.141785465198251" "0.3260061663372782"
Set-Row12Stats $ws 11 23 68 "-0.1657843576414015" "0.1297187088392248" "-0.228149078380763"  "0.1110368177920856"
Set-Row12Stats $ws 12 23 68 "0.2135201189924139"  "0.04951827186158275" "0.2585375793536707" "0.06985387930445178"

# ---------------------------------------------------------------------
# typestate_checker
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("typestate_checker")

Set-Row12Stats $ws 10 $null 520 "0.08993343971074573" "0.3939555332313259" "0.1261078291007821"  "0.3828516225378101"
Set-Row12Stats $ws 11 $null 520 "0.07727463833790227" "0.4435614692500816" "0.09527592158164427" "0.5104337983014946"
Set-Row12Stats $ws 12 $null 520 "0.06408768791953229" "0.5225032783206667" "0.08576989054836207" "0.5536912247004897"

# ---------------------------------------------------------------------
# infer
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("infer")

Set-Row12Stats $ws 10 2 2 "-0.07463933708620761" "0.547733910068501"  "-0.08587989564247843" "0.5531808807861933"
Set-Row12Stats $ws 11 2 2 "-0.02642855544759036" "0.8236209225496525" "-0.03184146471615851" "0.8262496889100787"
Set-Row12Stats $ws 12 2 2 "0.2391168558431198"   "0.04237662250330112" "0.2899717804431688"  "0.04108542769241475"

# col J width 20.7109375 -> 19.7109375, col K width 20.7109375 -> 21.7109375
$ws.Columns.Item(10).ColumnWidth = 18.75
$ws.Columns.Item(11).ColumnWidth = 20.75

# ---------------------------------------------------------------------
# openjml
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("openjml")

Set-Row12Stats $ws 10 $null 162 "-0.02930791451885823" "0.7910775440874591" "-0.03772642467808281" "0.7990497364707935"
Set-Row12Stats $ws 11 $null 162 "0.02099141080552319"  "0.8426468812248378" "0.02967862816832468"  "0.8412901317714291"
Set-Row12Stats $ws 12 $null 162 "0.2256752591064341"   "0.03184723740943501" "0.3154630979708622"  "0.02895586446936517"
